$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 197.5433703333333
$ws.Range("H2").Value = 592.6301109999999
$ws.Range("I2").Value = 0.3388703761585983
$ws.Range("J2").Value = 0.3388703761585982
$ws.Range("M2").Value = 0.01339666666666667
$ws.Range("N2").Value = 0.04019
$ws.Range("O2").Value = 0.08393217762128816
$ws.Range("P2").Value = 0.08393217762128814
$ws.Range("Q2").Value = 2.646422684565555
$ws.Range("R2").Value = 23.81780416109
$ws.Range("S2").Value = 0.0284421286023362
$ws.Range("T2").Value = 0.02844212860233619
$ws.Range("G3").Value = 197.5433703333333
$ws.Range("H3").Value = 592.6301109999999
$ws.Range("I3").Value = 0.3388703761585983
$ws.Range("J3").Value = 0.3388703761585982
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1260863333333333
$ws.Range("N3").Value = 0.378259
$ws.Range("O3").Value = 0.7899502755623498
$ws.Range("P3").Value = 0.7899502755623498
$ws.Range("Q3").Value = 24.90751923963877
$ws.Range("R3").Value = 224.167673156749
$ws.Range("S3").Value = 0.2676907470264018
$ws.Range("T3").Value = 0.2676907470264018
$ws.Range("G4").Value = 197.5433703333333
$ws.Range("H4").Value = 592.6301109999999
$ws.Range("I4").Value = 0.3388703761585983
$ws.Range("J4").Value = 0.3388703761585982
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02013
$ws.Range("N4").Value = 0.06039
$ws.Range("O4").Value = 0.1261175468163621
$ws.Range("P4").Value = 0.1261175468163621
$ws.Range("Q4").Value = 3.976548044809999
$ws.Range("R4").Value = 35.78893240329
$ws.Range("S4").Value = 0.04273750052986024
$ws.Range("T4").Value = 0.04273750052986023
$ws.Range("I5").Value = 0.1369374790620155
$ws.Range("J5").Value = 0.1369374790620154
$ws.Range("M5").Value = 0.01339666666666667
$ws.Range("N5").Value = 0.04019
$ws.Range("O5").Value = 0.08393217762128816
$ws.Range("P5").Value = 0.08393217762128814
$ws.Range("Q5").Value = 1.069419094891111
$ws.Range("R5").Value = 9.624771854019999
$ws.Range("S5").Value = 0.01149346081564451
$ws.Range("T5").Value = 0.01149346081564451
$ws.Range("I6").Value = 0.1369374790620155
$ws.Range("J6").Value = 0.1369374790620154
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1260863333333333
$ws.Range("N6").Value = 0.378259
$ws.Range("O6").Value = 0.7899502755623498
$ws.Range("P6").Value = 0.7899502755623498
$ws.Range("Q6").Value = 10.06512558881355
$ws.Range("R6").Value = 90.58613029932199
$ws.Range("S6").Value = 0.1081737993198526
$ws.Range("T6").Value = 0.1081737993198526
$ws.Range("I7").Value = 0.1369374790620155
$ws.Range("J7").Value = 0.1369374790620154
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02013
$ws.Range("N7").Value = 0.06039
$ws.Range("O7").Value = 0.1261175468163621
$ws.Range("P7").Value = 0.1261175468163621
$ws.Range("Q7").Value = 1.60692259618
$ws.Range("R7").Value = 14.46230336562
$ws.Range("S7").Value = 0.01727021892651833
$ws.Range("T7").Value = 0.01727021892651833
$ws.Range("G8").Value = 148.824417
$ws.Range("H8").Value = 446.473251
$ws.Range("I8").Value = 0.2552967790580629
$ws.Range("J8").Value = 0.2552967790580629
$ws.Range("M8").Value = 0.01339666666666667
$ws.Range("N8").Value = 0.04019
$ws.Range("O8").Value = 0.08393217762128816
$ws.Range("P8").Value = 0.08393217762128814
$ws.Range("Q8").Value = 1.99375110641
$ws.Range("R8").Value = 17.94375995769
$ws.Range("S8").Value = 0.02142761460604409
$ws.Range("T8").Value = 0.02142761460604409
$ws.Range("G9").Value = 148.824417
$ws.Range("H9").Value = 446.473251
$ws.Range("I9").Value = 0.2552967790580629
$ws.Range("J9").Value = 0.2552967790580629
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1260863333333333
$ws.Range("N9").Value = 0.378259
$ws.Range("O9").Value = 0.7899502755623498
$ws.Range("P9").Value = 0.7899502755623498
$ws.Range("Q9").Value = 18.764725050001
$ws.Range("R9").Value = 168.882525450009
$ws.Range("S9").Value = 0.2016717609670971
$ws.Range("T9").Value = 0.2016717609670971
$ws.Range("G10").Value = 148.824417
$ws.Range("H10").Value = 446.473251
$ws.Range("I10").Value = 0.2552967790580629
$ws.Range("J10").Value = 0.2552967790580629
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02013
$ws.Range("N10").Value = 0.06039
$ws.Range("O10").Value = 0.1261175468163621
$ws.Range("P10").Value = 0.1261175468163621
$ws.Range("Q10").Value = 2.99583551421
$ws.Range("R10").Value = 26.96251962789
$ws.Range("S10").Value = 0.03219740348492169
$ws.Range("T10").Value = 0.03219740348492169
$ws.Range("G11").Value = 35.426853
$ws.Range("H11").Value = 106.280559
$ws.Range("I11").Value = 0.06077202683121193
$ws.Range("J11").Value = 0.06077202683121192
$ws.Range("M11").Value = 0.01339666666666667
$ws.Range("N11").Value = 0.04019
$ws.Range("O11").Value = 0.08393217762128816
$ws.Range("P11").Value = 0.08393217762128814
$ws.Range("Q11").Value = 0.47460174069
$ws.Range("R11").Value = 4.271415666209999
$ws.Range("S11").Value = 0.005100728550402969
$ws.Range("T11").Value = 0.005100728550402968
$ws.Range("G12").Value = 35.426853
$ws.Range("H12").Value = 106.280559
$ws.Range("I12").Value = 0.06077202683121193
$ws.Range("J12").Value = 0.06077202683121192
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.1260863333333333
$ws.Range("N12").Value = 0.378259
$ws.Range("O12").Value = 0.7899502755623498
$ws.Range("P12").Value = 0.7899502755623498
$ws.Range("Q12").Value = 4.466841996309
$ws.Range("R12").Value = 40.201577966781
$ws.Range("S12").Value = 0.04800687934179838
$ws.Range("T12").Value = 0.04800687934179837
$ws.Range("G13").Value = 35.426853
$ws.Range("H13").Value = 106.280559
$ws.Range("I13").Value = 0.06077202683121193
$ws.Range("J13").Value = 0.06077202683121192
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.02013
$ws.Range("N13").Value = 0.06039
$ws.Range("O13").Value = 0.1261175468163621
$ws.Range("P13").Value = 0.1261175468163621
$ws.Range("Q13").Value = 0.71314255089
$ws.Range("R13").Value = 6.41828295801
$ws.Range("S13").Value = 0.007664418939010581
$ws.Range("T13").Value = 0.00766441893901058
$ws.Range("G14").Value = 121.3248153333333
$ws.Range("H14").Value = 363.974446
$ws.Range("I14").Value = 0.2081233388901116
$ws.Range("J14").Value = 0.2081233388901115
$ws.Range("M14").Value = 0.01339666666666667
$ws.Range("N14").Value = 0.04019
$ws.Range("O14").Value = 0.08393217762128816
$ws.Range("P14").Value = 0.08393217762128814
$ws.Range("Q14").Value = 1.625348109415556
$ws.Range("R14").Value = 14.62813298474
$ws.Range("S14").Value = 0.01746824504686039
$ws.Range("T14").Value = 0.01746824504686039
$ws.Range("G15").Value = 121.3248153333333
$ws.Range("H15").Value = 363.974446
$ws.Range("I15").Value = 0.2081233388901116
$ws.Range("J15").Value = 0.2081233388901115
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.1260863333333333
$ws.Range("N15").Value = 0.378259
$ws.Range("O15").Value = 0.7899502755623498
$ws.Range("P15").Value = 0.7899502755623498
$ws.Range("Q15").Value = 15.29740110772378
$ws.Range("R15").Value = 137.676609969514
$ws.Range("S15").Value = 0.1644070889071999
$ws.Range("T15").Value = 0.1644070889071999
$ws.Range("G16").Value = 121.3248153333333
$ws.Range("H16").Value = 363.974446
$ws.Range("I16").Value = 0.2081233388901116
$ws.Range("J16").Value = 0.2081233388901115
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02013
$ws.Range("N16").Value = 0.06039
$ws.Range("O16").Value = 0.1261175468163621
$ws.Range("P16").Value = 0.1261175468163621
$ws.Range("Q16").Value = 2.44226853266
$ws.Range("R16").Value = 21.98041679394
$ws.Range("S16").Value = 0.02624800493605123
$ws.Range("T16").Value = 0.02624800493605123
